$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26
$ws.Range("A26").Value = 111934989
$ws.Range("B26").Value = 77267
$ws.Range("C26").Value = "Ovaliderad"
$ws.Range("D26").Value = "NT"
$ws.Range("E26").Value = 6446
$ws.Range("F26").Value = "Kolflarnlav"
$ws.Range("G26").Value = "Carbonicola anthracophila"
$ws.Range("H26").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("P26").Value = "Skumsåstjärnen (Skumsåstjärnen), Ång"
$ws.Range("Q26").Value = 584600.4316807063
$ws.Range("R26").Value = 7048295.79915637
$ws.Range("S26").Value = 1
$ws.Range("T26").Value = "Västernorrland"
$ws.Range("U26").Value = "Sollefteå"
$ws.Range("V26").Value = "Ångermanland"
$ws.Range("W26").Value = "Ramsele"
# date-like text must stay text, not auto-convert to a date serial
$ws.Range("Y26").Value = "'2023-09-06"
$ws.Range("Z26").Value = "00:00"
$ws.Range("AA26").Value = "'2023-09-06"
$ws.Range("AB26").Value = "00:00"
$ws.Range("AD26").Value = $false
$ws.Range("AE26").Value = $false
$ws.Range("AG26").Value = $false
$ws.Range("AW26").Value = "Kamilla Andersson"
$ws.Range("AX26").Value = "Kamilla Andersson"

# Row 27
$ws.Range("A27").Value = 111935024
$ws.Range("B27").Value = 96348
$ws.Range("C27").Value = "Ovaliderad"
$ws.Range("D27").Value = "VU"
$ws.Range("E27").Value = 220787
$ws.Range("F27").Value = "Knärot"
$ws.Range("G27").Value = "Goodyera repens"
$ws.Range("H27").Value = "(L.) R. Br."
# numeric-looking text must stay text, not auto-convert to a number
$ws.Range("I27").Value = "'10"
$ws.Range("P27").Value = "Sollefteå (Sollefteå), Ång"
$ws.Range("Q27").Value = 584598.2684909205
$ws.Range("R27").Value = 7048259.615628711
$ws.Range("S27").Value = 2
$ws.Range("T27").Value = "Västernorrland"
$ws.Range("U27").Value = "Sollefteå"
$ws.Range("V27").Value = "Ångermanland"
$ws.Range("W27").Value = "Ramsele"
$ws.Range("Y27").Value = "'2023-09-06"
$ws.Range("Z27").Value = "19:51"
$ws.Range("AA27").Value = "'2023-09-06"
$ws.Range("AB27").Value = "19:51"
$ws.Range("AD27").Value = $false
$ws.Range("AE27").Value = $false
$ws.Range("AG27").Value = $false
$ws.Range("AW27").Value = "Kim Hultgren"
$ws.Range("AX27").Value = "Kim Hultgren"
